$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("A2").Value = 200
$ws.Range("D2").Value = 100

$ws.Range("A3").Value = 200
$ws.Range("D3").Value = 120

$ws.Range("D4").Value = 120

$ws.Range("E7").Value = "tekst"
